$wb = $excel.ActiveWorkbook

# 1. Insert a new "price" worksheet before "partOfPortfolio"
$before = $wb.Worksheets.Item("partOfPortfolio")
$price = $wb.Worksheets.Add($before)
$price.Name = "price"
$price.Range("A1").Value = "date"
$price.Range("B1").Value = "stockSymbol"
$price.Range("C1").Value = "stock price"
$price.Range("A2").Value = 0
$price.Range("B2").Value = 0

# 2. Update "partOfPortfolio" headers / sample row
$pop = $wb.Worksheets.Item("partOfPortfolio")
$pop.Range("B1").Value = "what part of the portfolio"
$pop.Range("A2").Value = 0
$pop.Range("B2").ClearContents() | Out-Null

# 3. Update "dowVSindex" sample row (date text -> numeric placeholder)
$dow = $wb.Worksheets.Item("dowVSindex")
$dow.Range("A2").Value = 0

# 4. Update "abserror" sample row (date text -> numeric placeholder)
$abs = $wb.Worksheets.Item("abserror")
$abs.Range("A2").Value = 0

# 5. Remove the "stockDataRep" worksheet entirely
$rep = $wb.Worksheets.Item("stockDataRep")
$rep.Delete() | Out-Null

# 6. Update "scalars_out" sample row value type (date text -> numeric placeholder)
$out = $wb.Worksheets.Item("scalars_out")
$out.Range("C2").Value = 0
